$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: issue number and report week date range ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Plain value updates (style/type unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 9
$ws.Range("K15").Value = -10
$ws.Range("L15").Value = 28.571428571428
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -78.571428571428
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = -28.571428571428
$ws.Range("L16").Value = -31.677018633540
$ws.Range("M16").Value = -50.226244343891
$ws.Range("N16").Value = -86.146095717884
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 15
$ws.Range("I17").Value = 236
$ws.Range("J17").Value = 279
$ws.Range("K17").Value = -15.412186379928
$ws.Range("L17").Value = 4.888888888888
$ws.Range("M17").Value = 43.902439024390
$ws.Range("N17").Value = -65.140324963072
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 86
$ws.Range("K18").Value = -9.302325581395
$ws.Range("L18").Value = -22
$ws.Range("M18").Value = -22.772277227722
$ws.Range("N18").Value = -90.382244143033
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -7.692307692307
$ws.Range("I19").Value = 386
$ws.Range("J19").Value = 468
$ws.Range("K19").Value = -17.521367521367
$ws.Range("L19").Value = 5.753424657534
$ws.Range("M19").Value = 32.646048109965
$ws.Range("N19").Value = -13.063063063063
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -90
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -54.237288135593
$ws.Range("N20").Value = -80.291970802919
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -29.166666666666
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -12.359550561797
$ws.Range("I21").Value = 851
$ws.Range("J21").Value = 1061
$ws.Range("K21").Value = -19.792648444863
$ws.Range("L21").Value = -8.396124865446
$ws.Range("M21").Value = 4.802955665024
$ws.Range("N21").Value = -71.024855294518
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = -13.636363636363
$ws.Range("M22").Value = 35.714285714285
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("H23").Value = -22.222222222222
$ws.Range("I23").Value = 58
$ws.Range("J23").Value = 76
$ws.Range("K23").Value = -23.684210526315
$ws.Range("L23").Value = -26.582278481012
$ws.Range("M23").Value = 31.818181818181
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 143.75
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -4.807692307692
$ws.Range("I24").Value = 1355
$ws.Range("J24").Value = 1191
$ws.Range("K24").Value = 13.769941225860
$ws.Range("L24").Value = 16.911130284728
$ws.Range("M24").Value = 45.230439442658
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 271.428571428571
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -13.333333333333
$ws.Range("I25").Value = 857
$ws.Range("J25").Value = 719
$ws.Range("K25").Value = 19.193324061196
$ws.Range("L25").Value = 28.101644245142
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 25.806451612903
$ws.Range("I26").Value = 392
$ws.Range("J26").Value = 413
$ws.Range("K26").Value = -5.084745762711
$ws.Range("L26").Value = 5.376344086021
$ws.Range("M26").Value = -17.299578059071
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = -31.25
$ws.Range("L27").Value = -15.384615384615
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = 3.125
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 11
$ws.Range("K29").Value = -42.105263157894
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -42.105263157894
$ws.Range("N29").Value = -88.888888888888
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = -35.714285714285
$ws.Range("L30").Value = -10
$ws.Range("M30").Value = -47.058823529411
$ws.Range("N30").Value = -89.285714285714

# --- Cells changing from numeric/percent to text placeholder ("0" or "***.*") style ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").PasteSpecial(-4122)

# --- Cells changing from text placeholder to numeric style ---
$ws.Range("C23").Value = 2
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false